$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells (row 1)
$ws.Range("F1").Value = "dob"
$ws.Range("G1").Value = "bloodGroup"

# New data cells (row 2)
# F2 mirrors D2: same "text date" formatting (numFmt 49 + quote-prefix),
# entered as quote-prefixed text so it stays a shared string, not a real date.
$ws.Range("F2").NumberFormat = $ws.Range("D2").NumberFormat
$ws.Range("F2").Value = "'2026-10-26"

$ws.Range("G2").Value = "B+"

# Remaining header cell (row 1)
$ws.Range("H1").Value = "tess"

# H2 uses the "text" number format (numFmt 49) but without the quote-prefix.
$ws.Range("H2").Value = "Testing Data"
$ws.Range("H2").NumberFormat = $ws.Range("D2").NumberFormat

# Move the active selection the way the recorded session ended up.
$ws.Range("H3").Select()
